$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Sending cluster ECs -> MuSCs, Target cluster stays ECs, TPM-derived metrics updated
$ws.Range("A2").Value = "MuSCs"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01868033333333333
$ws.Range("H2").Value = 0.056041
$ws.Range("M2").Value = 2.027115333333333
$ws.Range("N2").Value = 6.081346
$ws.Range("O2").Value = 0.006596284565418616
$ws.Range("P2").Value = 0.006596284565418615
$ws.Range("Q2").Value = 0.03786719013177778
$ws.Range("R2").Value = 0.340804711186
$ws.Range("S2").Value = 0.006596284565418616
$ws.Range("T2").Value = 0.006596284565418615

# Row 3: Sending cluster ECs -> MuSCs, Target cluster stays FAPs, TPM-derived metrics updated
$ws.Range("A3").Value = "MuSCs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01868033333333333
$ws.Range("H3").Value = 0.056041
$ws.Range("O3").Value = 0.8344762556643375
$ws.Range("P3").Value = 0.8344762556643374
$ws.Range("Q3").Value = 4.790465105061777
$ws.Range("R3").Value = 43.114185945556
$ws.Range("S3").Value = 0.8344762556643375
$ws.Range("T3").Value = 0.8344762556643374

# Row 4: Sending cluster ECs -> MuSCs, Target cluster stays MuSCs, TPM-derived metrics updated
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01868033333333333
$ws.Range("H4").Value = 0.056041
$ws.Range("M4").Value = 48.84026566666667
$ws.Range("N4").Value = 146.520797
$ws.Range("O4").Value = 0.158927459770244
$ws.Range("P4").Value = 0.158927459770244
$ws.Range("Q4").Value = 0.9123524427418891
$ws.Range("R4").Value = 8.211171984677001
$ws.Range("S4").Value = 0.158927459770244
$ws.Range("T4").Value = 0.158927459770244
